$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3665.55
$ws.Range("I70").Value = 4753.846
$ws.Range("J70").Value = 1644.4286
$ws.Range("K70").Value = 14261.538
$ws.Range("L70").Value = 4933.2858
$ws.Range("M70").Value = -13991.538
$ws.Range("N70").Value = -5473.2858
$ws.Range("H73").Value = 3665.55
$ws.Range("I73").Value = 4753.846
$ws.Range("J73").Value = 1644.4286
$ws.Range("K73").Value = 14261.538
$ws.Range("L73").Value = 4933.2858
$ws.Range("M73").Value = -13325.538
$ws.Range("N73").Value = -6805.2858
$ws.Range("H116").Value = 4223.25
$ws.Range("I116").Value = 4300
$ws.Range("J116").Value = 3993
$ws.Range("K116").Value = 4300
$ws.Range("L116").Value = 3993
$ws.Range("M116").Value = -858
$ws.Range("N116").Value = -10877
$ws.Range("H129").Value = 1006.37256
$ws.Range("I129").Value = 215
$ws.Range("J129").Value = 1055.8334
$ws.Range("K129").Value = 645
$ws.Range("L129").Value = 3167.5002
$ws.Range("M129").Value = 4355
$ws.Range("N129").Value = -13167.5002
$ws.Range("H141").Value = 2202.1538
$ws.Range("I141").Value = 2063.739
$ws.Range("J141").Value = 3263.3333
$ws.Range("K141").Value = 6191.217000000001
$ws.Range("L141").Value = 9789.999899999999
$ws.Range("M141").Value = -1011.217000000001
$ws.Range("N141").Value = -20149.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 915.8333
$ws.Range("I2").Value = 927.2857
$ws.Range("J2").Value = 875.75
$ws.Range("K2").Value = 927.2857
$ws.Range("L2").Value = 875.75
$ws.Range("M2").Value = -814.2857
$ws.Range("N2").Value = -1101.75
$ws.Range("H32").Value = 4434.8203
$ws.Range("I32").Value = 3128.3845
$ws.Range("J32").Value = 10967
$ws.Range("K32").Value = 3128.3845
$ws.Range("L32").Value = 10967
$ws.Range("M32").Value = -2841.3845
$ws.Range("N32").Value = -11541
$ws.Range("H97").Value = 750.8261
$ws.Range("I97").Value = 783.7143
$ws.Range("J97").Value = 405.5
$ws.Range("K97").Value = 783.7143
$ws.Range("L97").Value = 405.5
$ws.Range("M97").Value = -287.7143
$ws.Range("N97").Value = -1397.5
$ws.Range("H116").Value = 915.8333
$ws.Range("I116").Value = 927.2857
$ws.Range("J116").Value = 875.75
$ws.Range("K116").Value = 927.2857
$ws.Range("L116").Value = 875.75
$ws.Range("M116").Value = 1366.7143
$ws.Range("N116").Value = -5463.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 915.8333
$ws.Range("I3").Value = 927.2857
$ws.Range("J3").Value = 875.75
$ws.Range("K3").Value = 927.2857
$ws.Range("L3").Value = 875.75
$ws.Range("M3").Value = -813.2857
$ws.Range("N3").Value = -1103.75
$ws.Range("H94").Value = 674.125
$ws.Range("I94").Value = 548.10345
$ws.Range("J94").Value = 1006.36365
$ws.Range("K94").Value = 548.10345
$ws.Range("L94").Value = 1006.36365
$ws.Range("M94").Value = -97.10344999999995
$ws.Range("N94").Value = -1908.36365
$ws.Range("H99").Value = 2498.2334
$ws.Range("I99").Value = 2589.8928
$ws.Range("J99").Value = 1215
$ws.Range("K99").Value = 2589.8928
$ws.Range("L99").Value = 1215
$ws.Range("M99").Value = -1091.8928
$ws.Range("N99").Value = -4211

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 376.27274
$ws.Range("I22").Value = 317.66666
$ws.Range("K22").Value = 317.66666
$ws.Range("M22").Value = 32.33334000000002
$ws.Range("H31").Value = 3271.8215
$ws.Range("I31").Value = 1981.6818
$ws.Range("J31").Value = 8002.3335
$ws.Range("K31").Value = 1981.6818
$ws.Range("L31").Value = 8002.3335
$ws.Range("M31").Value = -1686.6818
$ws.Range("N31").Value = -8592.333500000001
$ws.Range("H34").Value = 3271.8215
$ws.Range("I34").Value = 1981.6818
$ws.Range("J34").Value = 8002.3335
$ws.Range("K34").Value = 1981.6818
$ws.Range("L34").Value = 8002.3335
$ws.Range("M34").Value = -1779.6818
$ws.Range("N34").Value = -8406.333500000001
$ws.Range("H88").Value = 10000
$ws.Range("J88").Value = 10000
$ws.Range("L88").Value = 10000
$ws.Range("N88").Value = -10812
$ws.Range("H91").Value = 10000
$ws.Range("J91").Value = 10000
$ws.Range("L91").Value = 10000
$ws.Range("N91").Value = -12808
$ws.Range("H94").Value = 4303.0415
$ws.Range("I94").Value = 852.8461
$ws.Range("J94").Value = 8380.546
$ws.Range("K94").Value = 852.8461
$ws.Range("L94").Value = 8380.546
$ws.Range("M94").Value = -401.8461
$ws.Range("N94").Value = -9282.546
$ws.Range("H122").Value = 1698.909
$ws.Range("I122").Value = 1636.125
$ws.Range("J122").Value = 1866.3334
$ws.Range("K122").Value = 4908.375
$ws.Range("L122").Value = 5599.0002
$ws.Range("M122").Value = -2458.375
$ws.Range("N122").Value = -10499.0002
$ws.Range("H134").Value = 1926.0358
$ws.Range("I134").Value = 1009.9
$ws.Range("J134").Value = 2435
$ws.Range("K134").Value = 3029.7
$ws.Range("L134").Value = 7305
$ws.Range("M134").Value = -494.6999999999998
$ws.Range("N134").Value = -12375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 51.565216
$ws.Range("J12").Value = 61.35294
$ws.Range("L12").Value = 184.05882
$ws.Range("N12").Value = -530.05882
$ws.Range("H63").Value = 3133.3333
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 3133.3333
$ws.Range("K63").Value = 0
$ws.Range("L63").ClearContents()
$ws.Range("M63").Value = 9399.999899999999
$ws.Range("N63").Value = -10897.9999
$ws.Range("H64").Value = 2152.8
$ws.Range("I64").Value = 1127.4286
$ws.Range("J64").Value = 3050
$ws.Range("K64").Value = 3382.2858
$ws.Range("L64").Value = 9150
$ws.Range("M64").Value = -3112.2858
$ws.Range("N64").Value = -9690
$ws.Range("H66").Value = 3133.3333
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 3133.3333
$ws.Range("K66").Value = 0
$ws.Range("L66").ClearContents()
$ws.Range("M66").Value = 28199.9997
$ws.Range("N66").Value = -35687.9997
$ws.Range("H67").Value = 2152.8
$ws.Range("I67").Value = 1127.4286
$ws.Range("J67").Value = 3050
$ws.Range("K67").Value = 3382.2858
$ws.Range("L67").Value = 9150
$ws.Range("M67").Value = -2446.2858
$ws.Range("N67").Value = -11022
$ws.Range("H107").Value = 845.25
$ws.Range("I107").Value = 262.25
$ws.Range("K107").Value = 786.75
$ws.Range("M107").Value = 1133.25
$ws.Range("H113").Value = 6579434
$ws.Range("J113").Value = 519.4286
$ws.Range("L113").Value = 1558.2858
$ws.Range("N113").Value = -5898.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 544.2857
$ws.Range("J97").Value = 643.3333
$ws.Range("L97").Value = 643.3333
$ws.Range("N97").Value = -1635.3333
$ws.Range("H126").Value = 2412.0908
$ws.Range("I126").Value = 2208.9
$ws.Range("J126").Value = 4444
$ws.Range("K126").Value = 6626.700000000001
$ws.Range("L126").Value = 13332
$ws.Range("M126").Value = -4156.700000000001
$ws.Range("N126").Value = -18272
$ws.Range("H137").Value = 28000
$ws.Range("J137").Value = 28000
$ws.Range("L137").Value = 28000
$ws.Range("N137").Value = -38200

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2552.3157
$ws.Range("I7").Value = 2261.5386
$ws.Range("J7").Value = 3182.3333
$ws.Range("K7").Value = 2261.5386
$ws.Range("L7").Value = 3182.3333
$ws.Range("M7").Value = -2149.5386
$ws.Range("N7").Value = -3406.3333
$ws.Range("H61").Value = 1597.6
$ws.Range("I61").Value = 1597.6
$ws.Range("K61").Value = 1597.6
$ws.Range("M61").Value = -1395.6
$ws.Range("H113").Value = 1597.6
$ws.Range("I113").Value = 1597.6
$ws.Range("K113").Value = 1597.6
$ws.Range("M113").Value = 572.4000000000001
$ws.Range("H126").Value = 2552.3157
$ws.Range("I126").Value = 2261.5386
$ws.Range("J126").Value = 3182.3333
$ws.Range("K126").Value = 6784.6158
$ws.Range("L126").Value = 9546.999899999999
$ws.Range("M126").Value = -4314.6158
$ws.Range("N126").Value = -14486.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 215.61539
$ws.Range("I113").Value = 215.61539
$ws.Range("K113").Value = 646.84617
$ws.Range("M113").Value = 1523.15383
$ws.Range("H132").Value = 2046.6765
$ws.Range("I132").Value = 929.45
$ws.Range("J132").Value = 3642.7144
$ws.Range("K132").Value = 2788.35
$ws.Range("L132").Value = 10928.1432
$ws.Range("M132").Value = -258.3500000000004
$ws.Range("N132").Value = -15988.1432
